$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424853956436465
$ws.Range("D2").Value = 0.001537061079751822
$ws.Range("E2").Value = 0.971244539797226
$ws.Range("F2").Value = 0.4802030090081786
$ws.Range("G2").Value = 0.3799583934371782
$ws.Range("H2").Value = 0.3921776620666293
$ws.Range("M2").Value = 5.46345062891146

$ws.Range("B3").Value = 0.1329189329742348
$ws.Range("D3").Value = 0.001190009458271746
$ws.Range("E3").Value = 0.8535059057366254
$ws.Range("F3").Value = 0.4573835664836992
$ws.Range("G3").Value = 0.3512005984779165
$ws.Range("H3").Value = 0.3913316344414568
$ws.Range("M3").Value = 4.771688397558961

$ws.Range("B4").Value = 0.1271164234097313
$ws.Range("D4").Value = 0.001006887159459069
$ws.Range("E4").Value = 0.7814208774350959
$ws.Range("F4").Value = 0.4447363980145411
$ws.Range("G4").Value = 0.3348924932047197
$ws.Range("H4").Value = 0.391826831173006
$ws.Range("M4").Value = 4.346002943016572

$ws.Range("B5").Value = 0.1247698552961083
$ws.Range("D5").Value = 0.0009393067780596454
$ws.Range("E5").Value = 0.7520931493393874
$ws.Range("F5").Value = 0.4399148614288038
$ws.Range("G5").Value = 0.3285729282480077
$ws.Range("H5").Value = 0.3922777851151977
$ws.Range("M5").Value = 4.172274232350247

$ws.Range("B6").Value = 0.124381298706794
$ws.Range("D6").Value = 0.0009284963697879789
$ws.Range("E6").Value = 0.747226048879952
$ws.Range("F6").Value = 0.4391340077553139
$ws.Range("G6").Value = 0.3275428846087891
$ws.Range("H6").Value = 0.3923675406459495
$ws.Range("M6").Value = 4.143410506719647

$ws.Range("B7").Value = 0.1270847037940825
$ws.Range("D7").Value = 0.001005947869197499
$ws.Range("E7").Value = 0.7810251667165744
$ws.Range("F7").Value = 0.4446700417160798
$ws.Range("G7").Value = 0.3348059621580433
$ws.Range("H7").Value = 0.3918319119144087
$ws.Range("M7").Value = 4.343661050602037

$ws.Range("B8").Value = 0.1391721058803341
$ws.Range("D8").Value = 0.001410877104804698
$ws.Range("E8").Value = 0.9306023781991541
$ws.Range("F8").Value = 0.4720452510194519
$ws.Range("G8").Value = 0.3697546665290332
$ws.Range("H8").Value = 0.3916717280169877
$ws.Range("M8").Value = 5.225112178061124

$ws.Range("B9").Value = 0.1634404375874112
$ws.Range("D9").Value = 0.002465542594332959
$ws.Range("E9").Value = 1.225810626037912
$ws.Range("F9").Value = 0.5370431516352596
$ws.Range("G9").Value = 0.4495992380827829
$ws.Range("H9").Value = 0.3996824639121144
$ws.Range("M9").Value = 6.94735415644368

$ws.Range("B10").Value = 0.1816159560377599
$ws.Range("D10").Value = 0.003431821204104679
$ws.Range("E10").Value = 1.444251710462851
$ws.Range("F10").Value = 0.592390378635713
$ws.Range("G10").Value = 0.5160092524539266
$ws.Range("H10").Value = 0.4110289445704325
$ws.Range("M10").Value = 8.21078158333313

$ws.Range("B11").Value = 0.1899600438936062
$ws.Range("D11").Value = 0.003920366556256383
$ws.Range("E11").Value = 1.544063417558505
$ws.Range("F11").Value = 0.6193728945512618
$ws.Range("G11").Value = 0.5480933860408186
$ws.Range("H11").Value = 0.4174624682584067
$ws.Range("M11").Value = 8.78561072456057

$ws.Range("B12").Value = 0.1931306643634798
$ws.Range("D12").Value = 0.004113047097506595
$ws.Range("E12").Value = 1.581931715700478
$ws.Range("F12").Value = 0.6298630025393237
$ws.Range("G12").Value = 0.5605282666423363
$ws.Range("H12").Value = 0.4200888349962781
$ws.Range("M12").Value = 9.003336873078808

$ws.Range("B13").Value = 0.192447330266404
$ws.Range("D13").Value = 0.004071197979117258
$ws.Range("E13").Value = 1.573772769200701
$ws.Range("F13").Value = 0.6275914454226239
$ws.Range("G13").Value = 0.5578372460689991
$ws.Range("H13").Value = 0.4195146266550012
$ws.Range("M13").Value = 8.956442753830856

$ws.Range("B14").Value = 0.19022067450679
$ws.Range("D14").Value = 0.003936060870993074
$ws.Range("E14").Value = 1.547177377055448
$ws.Range("F14").Value = 0.620230387387096
$ws.Range("G14").Value = 0.5491106006516304
$ws.Range("H14").Value = 0.4176746884441229
$ws.Range("M14").Value = 8.803521938790368

$ws.Range("B15").Value = 0.1888582018319482
$ws.Range("D15").Value = 0.003854304816998777
$ws.Range("E15").Value = 1.530896535030735
$ws.Range("F15").Value = 0.6157573916799919
$ws.Range("G15").Value = 0.5438029127341508
$ws.Range("H15").Value = 0.4165726515946062
$ws.Range("M15").Value = 8.709861321944629

$ws.Range("B16").Value = 0.1810721761904404
$ws.Range("D16").Value = 0.003400933004639839
$ws.Range("E16").Value = 1.437738358810265
$ws.Range("F16").Value = 0.5906644186344465
$ws.Range("G16").Value = 0.5139515316785719
$ws.Range("H16").Value = 0.4106347068877483
$ws.Range("M16").Value = 8.17322012435784

$ws.Range("B17").Value = 0.1763151336130164
$ws.Range("D17").Value = 0.00313579368659056
$ws.Range("E17").Value = 1.380707625208032
$ws.Range("F17").Value = 0.5757420045841712
$ws.Range("G17").Value = 0.4961297692935602
$ws.Range("H17").Value = 0.4073228285533332
$ws.Range("M17").Value = 7.844057123830339

$ws.Range("B18").Value = 0.1735861586969634
$ws.Range("D18").Value = 0.002987863333395069
$ws.Range("E18").Value = 1.347946099938866
$ws.Range("F18").Value = 0.5673281106473667
$ws.Range("G18").Value = 0.4860545686578632
$ws.Range("H18").Value = 0.4055372024969586
$ws.Range("M18").Value = 7.654737255696489

$ws.Range("B19").Value = 0.1726634036713506
$ws.Range("D19").Value = 0.002938542723104121
$ws.Range("E19").Value = 1.33686042007443
$ws.Range("F19").Value = 0.5645079638823063
$ws.Range("G19").Value = 0.4826729288121783
$ws.Range("H19").Value = 0.4049528893955596
$ws.Range("M19").Value = 7.590636723058537

$ws.Range("B20").Value = 0.1768207891223028
$ws.Range("D20").Value = 0.003163540636174744
$ws.Range("E20").Value = 1.386774338747898
$ws.Range("F20").Value = 0.5773129292181665
$ws.Range("G20").Value = 0.4980086590239239
$ws.Range("H20").Value = 0.407662987031955
$ws.Range("M20").Value = 7.879096202803453

$ws.Range("B21").Value = 0.1908744017464556
$ws.Range("D21").Value = 0.003975540318810644
$ws.Range("E21").Value = 1.554987071328583
$ws.Range("F21").Value = 0.6223850108650026
$ws.Range("G21").Value = 0.5516659560802282
$ws.Range("H21").Value = 0.4182099042378979
$ws.Range("M21").Value = 8.848436792042548

$ws.Range("B22").Value = 0.2001227509912269
$ws.Range("D22").Value = 0.004551301348367787
$ws.Range("E22").Value = 1.665346490277898
$ws.Range("F22").Value = 0.6534362974816332
$ws.Range("G22").Value = 0.5884048102143993
$ws.Range("H22").Value = 0.426214648476531
$ws.Range("M22").Value = 9.482272983370081

$ws.Range("B23").Value = 0.1951809291640672
$ws.Range("D23").Value = 0.004239666131097053
$ws.Range("E23").Value = 1.606404005937264
$ws.Range("F23").Value = 0.636713433369593
$ws.Range("G23").Value = 0.5686383704714899
$ws.Range("H23").Value = 0.4218382040413076
$ws.Range("M23").Value = 9.14394041958667

$ws.Range("B24").Value = 0.1765921638473884
$ws.Range("D24").Value = 0.003150982265909263
$ws.Range("E24").Value = 1.384031496554229
$ws.Range("F24").Value = 0.5766022005992539
$ws.Range("G24").Value = 0.4971586822676102
$ws.Range("H24").Value = 0.4075088330859558
$ws.Range("M24").Value = 7.863255286198125

$ws.Range("B25").Value = 0.1568146364597709
$ws.Range("D25").Value = 0.002149039155971977
$ws.Range("E25").Value = 1.145707788164344
$ws.Range("F25").Value = 0.5181772474600876
$ws.Range("G25").Value = 0.4267013635422927
$ws.Range("H25").Value = 0.3922777851151977
$ws.Range("M25").Value = 6.481925672891293
